$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("STRINGS_10")

# Row 156: "Journeymen Description" -> "Artisan Guilds Description"
$ws1.Range("B156").Value = "Artisan Guilds Description"
$ws1.Range("D156").Value = "The Artisan Guilds are a collective of master craftsmen and masons. Their presence in your kingdom lowers the costs of new construction, tools and training as well as attracting more like-minded individuals to settlements where they operate. Within a Masonry Guild, they can outfit and train capable Pioneers to expand the frontiers of the kingdom."

# Row 157: "ObjectData ProperName" / "Journeyman" / "PIONEER.INI ProperName" -> "Artisan Guilds ProperName" / "Artisan Guilds" / "Technology ProperName"
$ws1.Range("B157").Value = "Artisan Guilds ProperName"
$ws1.Range("D157").Value = "Artisan Guilds"
$ws1.Range("G157").Value = "Technology ProperName"

# Row 158: "Journeymen are expert craftsmen..." -> "Pioneers are elite settlers..." (now wraps to more lines)
$ws1.Range("D158").Value = "Pioneers are elite settlers, better trained and equipped to expand the kingdom’s borders. They are intrepid navigators in difficult terrain and take great pride in their work."
$ws1.Range("D158").WrapText = $true

# Row heights recalculated to fit the new, longer wrapped text
$ws1.Rows.Item(156).RowHeight = 64.9
$ws1.Rows.Item(157).RowHeight = 26.85
$ws1.Rows.Item(158).RowHeight = 39.55
